# Fill in "Values" (column E) definitions for many variables in Table1 / Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that get the common "0 = No; 1 = Yes; 99 = Unknown" value definition.
$rowsYesNoUnknown = @(
    111,112,113,114,115,116,117,118,119,120,121,122,123,124,
    126,128,129,130,131,
    216,217,218,219,220,221,224,226,227,
    256,257,258,259,260,261,262
)

# Row 163 (D02, gender recode) gets its own distinct value definition.
# Entered first so it lands at the lower new shared-string index.
$ws.Range("E163").Value = "Female; Male"
$ws.Rows(163).RowHeight = 16

# Row 125 (Comp08) gets its own distinct value definition.
$ws.Range("E125").Value = "0 = Default; 1 = Yes"
$ws.Rows(125).RowHeight = 16

foreach ($r in $rowsYesNoUnknown) {
    $ws.Range("E$r").Value = "0 = No; 1 = Yes; 99 = Unknown"
    $ws.Rows($r).RowHeight = 16
}

# Reflect the view state captured in the saved file (active cell).
$ws.Range("E126").Select()

$wb.Save()
